# Apply the "cryptos list" refresh (GitHub Actions commit, Thu Jun 13 12:24:48 UTC 2024).
# Updates Coin/Link/Price/Volume(1h) cells in rows 2-51 of the active sheet to match
# the latest scrape. Price-column values that are purely numeric-looking text (e.g. "1.00",
# "0.0900") are forced to the Text number format first so Excel does not silently
# normalize them into floating point numbers and drop the original formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.786.71"
$ws.Range("E2").Value = "  -0.09%  "

# Row 3
$ws.Range("D3").Value = "3.499.59"
$ws.Range("E3").Value = "  -1.07%  "

# Row 4
$ws.Range("E4").Value = "  +0.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.92"
$ws.Range("E5").Value = "  -1.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.32"
$ws.Range("E6").Value = "  -1.34%  "

# Row 7
$ws.Range("D7").Value = "3.497.45"
$ws.Range("E7").Value = "  -1.14%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  +0.34%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.144"
$ws.Range("E10").Value = "  +2.93%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.55"
$ws.Range("E11").Value = "  +6.33%  "

# Row 12
$ws.Range("E12").Value = "  +0.80%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000215"
$ws.Range("E13").Value = "  -1.82%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.93"
$ws.Range("E14").Value = "  -0.23%  "

# Row 15
$ws.Range("D15").Value = "4.090.56"
$ws.Range("E15").Value = "  -1.05%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "67.716.99"
$ws.Range("E16").Value = "  +0.08%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.468.17"
$ws.Range("E17").Value = "  -2.08%  "

# Row 18
$ws.Range("E18").Value = "  -0.04%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.49"

# Row 20
$ws.Range("E20").Value = "  +0.47%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.96"
$ws.Range("E21").Value = "  +1.98%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "444.89"
$ws.Range("E22").Value = "  -0.58%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.624"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.90"
$ws.Range("E24").Value = "  +2.41%  "

# Row 25
$ws.Range("D25").Value = "3.639.28"
$ws.Range("E25").Value = "  -1.03%  "

# Row 26
$ws.Range("E26").Value = "  -0.03%  "

# Row 27
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000124"
$ws.Range("E27").Value = "  -4.79%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.70"
$ws.Range("E28").Value = "  +0.43%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.86"
$ws.Range("E29").Value = "  -3.29%  "

# Row 30
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.65"
$ws.Range("E30").Value = "  +2.24%  "

# Row 31
$ws.Range("E31").Value = "  -1.72%  "

# Row 32
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.169"
$ws.Range("E32").Value = "  -0.10%  "

# Row 33
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.01%  "

# Row 34
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.20"
$ws.Range("E34").Value = "  -0.12%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.54"
$ws.Range("E35").Value = "  -0.82%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.85"
$ws.Range("E36").Value = "  +0.21%  "

# Row 37
$ws.Range("D37").Value = "3.493.35"
$ws.Range("E37").Value = "  -0.90%  "

# Row 38
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.98"
$ws.Range("E38").Value = "  -1.11%  "

# Row 39
$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.01%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.32"
$ws.Range("E40").Value = "  +5.75%  "

# Row 41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.16%  "

# Row 42
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "176.25"
$ws.Range("E42").Value = "  +0.96%  "

# Row 43
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0900"
$ws.Range("E43").Value = "  +0.10%  "

# Row 44
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.41"
$ws.Range("E44").Value = "  +0.10%  "

# Row 45
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.896"
$ws.Range("E45").Value = "  +0.87%  "

# Row 46
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "30.08"
$ws.Range("E46").Value = "  +6.58%  "

# Row 47
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.87"
$ws.Range("E47").Value = "  +2.86%  "

# Row 48
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.28"
$ws.Range("E48").Value = "  -0.17%  "

# Row 49
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.53"
$ws.Range("E49").Value = "  -5.55%  "

# Row 50
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.59"
$ws.Range("E50").Value = "  -0.22%  "

# Row 51
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.991"
$ws.Range("E51").Value = "  -0.02%  "
